$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.5696291547463442
$ws.Range("C4").Value = 0.5640000000000001
$ws.Range("D4").Value = 0.6289259988206177
$ws.Range("E4").Value = 0.6075
$ws.Range("F4").Value = 0.6641157913352068
$ws.Range("G4").Value = 0.985
$ws.Range("H4").Value = 0.500948808950948
$ws.Range("I4").Value = 0.502
$ws.Range("J4").Value = 0.6091506016182076
$ws.Range("K4").Value = 0.65
$ws.Range("L4").Value = 0.5865631528644721
$ws.Range("M4").Value = 0.6110000000000001
$ws.Range("B5").Value = 0.6889796196133278
$ws.Range("C5").Value = 0.7180000000000001
$ws.Range("D5").Value = 0.674428656652411
$ws.Range("E5").Value = 0.6855
$ws.Range("F5").Value = 0.6473491909008886
$ws.Range("G5").Value = 0.9480000000000001
$ws.Range("H5").Value = 0.4918555057391506
$ws.Range("I5").Value = 0.487
$ws.Range("J5").Value = 0.6017323807093553
$ws.Range("K5").Value = 0.611
$ws.Range("L5").Value = 0.6032489817106115
$ws.Range("M5").Value = 0.6214999999999999
$ws.Range("B6").Value = 0.630258001061337
$ws.Range("C6").Value = 0.723
$ws.Range("D6").Value = 0.5801521608554855
$ws.Range("E6").Value = 0.593
$ws.Range("F6").Value = 0.663189793483707
$ws.Range("G6").Value = 0.9860000000000001
$ws.Range("H6").Value = 0.4996475575957275
$ws.Range("I6").Value = 0.4995
$ws.Range("J6").Value = 0.6302343615535737
$ws.Range("K6").Value = 0.6380000000000001
$ws.Range("L6").Value = 0.6362502519737273
$ws.Range("M6").Value = 0.6475
$ws.Range("B7").Value = 0.2154210467690116
$ws.Range("C7").Value = 0.176
$ws.Range("D7").Value = 0.3613725041177805
$ws.Range("E7").Value = 0.421
$ws.Range("F7").Value = 0.07849586080043762
$ws.Range("G7").Value = 0.102
$ws.Range("H7").Value = 0.2328125
$ws.Range("I7").Value = 0.4865
$ws.Range("J7").Value = 0.4032295676642556
$ws.Range("K7").Value = 0.4069999999999999
$ws.Range("L7").Value = 0.4033372780918357
$ws.Range("M7").Value = 0.4135
